$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$secret = "U2FsdGVkX1+s1oRPt0B7OcjpVP5f3IBhIA53DP6hmaIufy/vOcfu72zH6t1I7wrJH9m3kgsE3/e9RQkIhB+/5g=="

$times = @("2:55-3:0", "3:0-3:5", "10:55-11:0", "11:0-11:5", "12:25-12:30", "12:50-12:55", "14:10-14:15", "14:15-14:20", "21:25-21:30", "21:30-21:35")

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Range("C$row").Value = $times[$i]
}

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $secret
}

$ws.Range("B12").Select()
